$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.373.93'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').Value = '1.844.20'
$ws.Range('E3').Value = '  -0.31%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9975'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.33%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.41'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6321'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.68%  '
$ws.Range('E7').Value = '  -0.26%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07500'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.87%  '
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.43'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.33%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07717'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.45%  '
$ws.Range('D12').Value = '1.844.15'
$ws.Range('E12').Value = '  -2.33%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.004'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6797'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001030'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -3.75%  '
$ws.Range('E16').Value = '  -1.40%  '
$ws.Range('D17').Value = '2.104.92'
$ws.Range('E17').Value = '  -3.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.159'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').Value = '29.393.33'
$ws.Range('E19').Value = '  -0.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '229.33'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.34'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9989'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.458'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9993'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '158.85'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.59%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.415'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('E28').Value = '  -0.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.06398'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +14.45%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.385'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.474'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.62%  '
$ws.Range('E32').Value = '  -0.88%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.046'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.39%  '
$ws.Range('E34').Value = '  -1.07%  '
$ws.Range('E35').Value = '  -1.92%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6993'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.55%  '
$ws.Range('E37').Value = '  -0.42%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.836'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +4.06%  '
$ws.Range('D39').Value = '1.256.00'
$ws.Range('E39').Value = '  +1.70%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01820'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.592'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.74%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9077'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.35%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9982'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.36%  '
$ws.Range('D44').Value = '2.006.23'
$ws.Range('E44').Value = '  -18.47%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.38'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.33%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '66.32'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.37%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000119'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.71%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1183'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +3.22%  '
$ws.Range('B49').Value = 'Aptos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.044'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.84%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.702'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.35%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.045'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.68%  '
